$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Fix 1: "if the execution fail." -> "if the execution fails." ---
$old1 = "if the execution fail."
$new1 = "if the execution fails."

# --- Fix 2: "Specify asynchronously delegate for timeout notification " -> "Specify asynchronous delegate for notification " ---
$old2 = "Specify asynchronously delegate for timeout notification "
$new2 = "Specify asynchronous delegate for notification "

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -ne -1) {
        continue
    }
    $tf = $sh.TextFrame
    $tr = $tf.TextRange
    $full = $tr.Text

    if ($full.Contains($old1)) {
        $idx = $full.IndexOf($old1)
        $chars = $tr.Characters($idx + 1, $old1.Length)
        $chars.Text = $new1
    }

    if ($full.Contains($old2)) {
        $idx = $full.IndexOf($old2)
        $chars = $tr.Characters($idx + 1, $old2.Length)
        $chars.Text = $new2
    }
}
